$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose data relocates to a different row (pure permutation; all other rows keep their
# own data untouched). Values come from the authoritative before/after diff.
$moveTargets = @(6,7,10,11,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,45,46,47,48,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82)

# before-row (value) that supplies the data for each after-row (key).
$map = @{}
$map[6] = 7
$map[7] = 6
$map[10] = 11
$map[11] = 10
$map[27] = 66
$map[28] = 67
$map[29] = 62
$map[30] = 56
$map[31] = 55
$map[32] = 58
$map[33] = 79
$map[34] = 78
$map[35] = 29
$map[36] = 27
$map[37] = 71
$map[38] = 28
$map[39] = 30
$map[40] = 80
$map[41] = 63
$map[42] = 52
$map[43] = 68
$map[45] = 82
$map[46] = 81
$map[47] = 65
$map[48] = 31
$map[50] = 48
$map[51] = 72
$map[52] = 33
$map[53] = 32
$map[54] = 59
$map[55] = 37
$map[56] = 57
$map[57] = 36
$map[58] = 35
$map[59] = 45
$map[60] = 39
$map[61] = 46
$map[62] = 64
$map[63] = 43
$map[64] = 60
$map[65] = 69
$map[66] = 73
$map[67] = 41
$map[68] = 38
$map[69] = 54
$map[70] = 50
$map[71] = 34
$map[72] = 53
$map[73] = 61
$map[74] = 42
$map[75] = 70
$map[76] = 74
$map[77] = 47
$map[78] = 51
$map[79] = 40
$map[80] = 75
$map[81] = 76
$map[82] = 77

# Snapshot the source rows (A:Z) BEFORE any writes, since several rows both donate and
# receive data (e.g. rows 6 & 7 swap).
$snap = @{}
foreach ($r in $moveTargets) {
    $src = $map[$r]
    if (-not $snap.ContainsKey($src)) {
        $snap[$src] = $ws.Range("A$src`:Z$src").Formula
    }
}

# Write each row to its new position.
foreach ($r in $moveTargets) {
    $src = $map[$r]
    $ws.Range("A$r`:Z$r").Formula = $snap[$src]
}

# "Förändrad" (C) column: every data row (2-82) advances from 46059 to 46060.
for ($r = 2; $r -le 82; $r++) {
    $ws.Range("C$r").Value = 46060
}
